$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_11_2")

$ws.Range("A22").Value = 'Main air reservoir△1'
$ws.Range("B22:G22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("H22").Value = '(+)reservoir & (-)reservoir△1'
$ws.Range("I22:L22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("M22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("N22:AC22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AD22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AE22:AI22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AJ22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AK22:AL22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AM22").Value = '(+)condenser & (+)reservoir & (+). & (-)reservoir△1'
$ws.Range("AN22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AO22").Value = '(+)condenser & (+)reservoir & (+). & (-)reservoir△1'
$ws.Range("AP22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AQ22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AR22").NumberFormat = "@"
$ws.Range("AR22").Value = '1'
$ws.Range("AR22").ClearFormats()
$ws.Range("AS22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("AT22:AU22").Value = '(+)condenser & (+)reservoir & (+). & (-)reservoir△1'
$ws.Range("AV22:BZ22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CA22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CB22:CG22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CH22").Value = '(+)11.2 & (+)PAINTING & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CI22:CQ22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CR22").Value = '(+)drain & (+)cooler & (+)△2 & (-)Main & (-)air & (-)reservoir△1'
$ws.Range("CS22:CX22").Value = '(+)Maker & (+)’ & (+)standard & (+). & (-)Main & (-)air & (-)reservoir△1'